$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D so the existing "Tipo" column shifts to E
# (the new D column inherits the bold/bordered header formatting)
$ws.Columns("D").Insert()

# New header and data for the inserted MAE column
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.2529908105885643

# Updated MSE and R2 values
$ws.Range("B2").Value = 0.09766946137408543
$ws.Range("C2").Value = 0.9986501815983999
